# Updates cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.863.71'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '2.579.91'
$ws.Range('E3').Value = '  +1.12%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.62'
$ws.Range('E5').Value = '  +0.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.52'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.595'
$ws.Range('E8').Value = '  +2.18%  '
$ws.Range('E9').Value = '  +2.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.65'
$ws.Range('E10').Value = '  +2.53%  '
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.352'
$ws.Range('E12').Value = '  -0.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.22'
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').Value = '3.044.88'
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('D15').Value = '62.780.92'
$ws.Range('E15').Value = '  -0.30%  '
$ws.Range('E16').Value = '  +2.54%  '
$ws.Range('D17').Value = '2.576.69'
$ws.Range('E17').Value = '  +1.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.29'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '341.71'
$ws.Range('E19').Value = '  +1.76%  '
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.66'
$ws.Range('E21').Value = '  -1.31%  '
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('E23').Value = '  -1.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.06'
$ws.Range('E24').Value = '  +2.58%  '
$ws.Range('D25').Value = '2.712.26'
$ws.Range('E25').Value = '  +1.14%  '
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.58'
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('E29').Value = '  +6.65%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.30'
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.43'
$ws.Range('E31').Value = '  -3.54%  '
$ws.Range('E32').Value = '  +1.94%  '
$ws.Range('E33').Value = '  +0.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '463.38'
$ws.Range('E34').Value = '  +13.95%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '175.05'
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('E36').Value = '  +3.41%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.397'
$ws.Range('E38').Value = '  -0.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.00'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.50'
$ws.Range('E40').Value = '  +3.88%  '
$ws.Range('E42').Value = '  -2.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '157.89'
$ws.Range('E43').Value = '  +4.40%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.636'
$ws.Range('E45').Value = '  +5.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.08'
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0965'
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('E49').Value = '  -0.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.41'
$ws.Range('E50').Value = '  +0.65%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.45'
$ws.Range('E51').Value = '  +1.31%  '
